# Update: Threat Alert Report - 2026-01-25 09:07
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date cells (A column) ---
# Assign as text (leading apostrophe forces text so Excel does not
# auto-convert these "dd-MMM-yy" looking strings into real date serials).
$ws.Range("A2").Value = "'30-JAN-26"
$ws.Range("A3").Value = "'20-FEB-26"
$ws.Range("A4").Value = "'27-FEB-26"
$ws.Range("A6").Value = "'27-MAR-26"
$ws.Range("A7").Value = "'29-MAR-26"
$ws.Range("A8").Value = "'02-APR-26"
$ws.Range("A9").Value = "'17-MAY-26"

# Restore the original (General / centered / bordered, non quote-prefixed)
# cell formatting on those date cells by copying it from an untouched
# date cell in the same column (A5 keeps its original formatting).
$ws.Range("A5").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)

# --- Row 2 ---
$ws.Range("D2").Value = 462
$ws.Range("E2").Value = 475
$ws.Range("F2").Value = -13

# --- Row 3 ---
$ws.Range("D3").Value = 462
$ws.Range("F3").Value = -57

# --- Row 4 ---
$ws.Range("D4").Value = 462
$ws.Range("E4").Value = 519
$ws.Range("F4").Value = -57

# --- Row 5 ---
$ws.Range("D5").Value = 1240
$ws.Range("F5").Value = -218

# --- Row 6 ---
$ws.Range("D6").Value = 462
$ws.Range("E6").Value = 713
$ws.Range("F6").Value = -251
# IMPACT changes from LOW THREAT to MEDIUM THREAT - MONITOR: copy formatting from row 5 (MEDIUM)
$ws.Range("J5").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = "MEDIUM THREAT - MONITOR"

# --- Row 7 ---
$ws.Range("D7").Value = 462
$ws.Range("E7").Value = 621
$ws.Range("F7").Value = -159

# --- Row 8 ---
$ws.Range("D8").Value = 456
$ws.Range("E8").Value = 519
$ws.Range("F8").Value = -63

# --- Row 9 ---
$ws.Range("D9").Value = 960
$ws.Range("E9").Value = 1774
$ws.Range("F9").Value = -814
# IMPACT changes from LOW THREAT to HIGH THREAT ALERT - NEED ACTION: copy formatting from row 10 (HIGH)
$ws.Range("J10").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = "HIGH THREAT ALERT - NEED ACTION"

# --- Row 10 ---
$ws.Range("D10").Value = 1234
$ws.Range("F10").Value = -540

# --- Row 11 ---
$ws.Range("D11").Value = 1240
$ws.Range("F11").Value = -534
